# [UPDATE] Se actualizaron scripts y objetos para la ejecucion de las ENQ
# Adds new user/sucursal records to the "Users" sheet (rows 34-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# --- Row 34 : CCUENCA / 001 ---
$ws.Range("A34").Value = "CCUENCA"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").HorizontalAlignment = -4152
$ws.Range("C34").Value = "001"

# --- Row 35 : F04033 / 533 ---
$ws.Range("A35").Value = "F04033"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").HorizontalAlignment = -4152
$ws.Range("C35").Value = "533"

# --- Row 36 : F04169 / 369 ---
$ws.Range("A36").Value = "F04169"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").HorizontalAlignment = -4152
$ws.Range("C36").Value = "369"

# --- Row 37 : F00463 / 063 ---
$ws.Range("A37").Value = "F00463"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").HorizontalAlignment = -4152
$ws.Range("C37").Value = "063"

# --- Row 38 : JANDINO / 102 (Sucursal value entered last, see below) ---
$ws.Range("A38").Value = "JANDINO"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").HorizontalAlignment = -4152

# --- Row 39 : F00219 / 019 ---
$ws.Range("A39").Value = "F00219"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").HorizontalAlignment = -4152
$ws.Range("C39").Value = "019"

# --- Row 40 : F00089 / 089 ---
$ws.Range("A40").Value = "F00089"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").HorizontalAlignment = -4152
$ws.Range("C40").Value = "089"

# --- Rows 41-43 : blank placeholder rows, only formatted (text, right aligned) ---
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").HorizontalAlignment = -4152

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").HorizontalAlignment = -4152

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").HorizontalAlignment = -4152

# Fill in row 38's Sucursal value last so the new shared-string entry
# ends up appended at the end, matching the authored edit order.
$ws.Range("C38").Value = "102"

# Leave the sheet on the newly added record, mirroring the author's final
# cursor position/selection.
$ws.Activate()
$ws.Range("C34").Select()
